$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F13").Value = 'PAOK B'
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 'Kampaniakos'
$ws.Range("I13").Value = 3
$ws.Range("J13").Value = 1.46
$ws.Range("K13").Value = '30/09/2023 03:12'
$ws.Range("L13").Value = 1.52
$ws.Range("M13").Value = '30/09/2023 13:01'
$ws.Range("N13").Value = 4.27
$ws.Range("O13").Value = '30/09/2023 03:12'
$ws.Range("P13").Value = 4.05
$ws.Range("Q13").Value = '30/09/2023 13:17'
$ws.Range("R13").Value = 5.93
$ws.Range("S13").Value = '30/09/2023 03:12'
$ws.Range("T13").Value = 6.09
$ws.Range("U13").Value = '30/09/2023 13:17'
$ws.Range("V13").Value = 'https://www.betexplorer.com/football/greece/super-league-2/paok-kampaniakos/S6W0xkkj/'

$ws.Range("F14").Value = 'Panachaiki'
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 'Panathinaikos B'
$ws.Range("I14").Value = 2
$ws.Range("J14").Value = 2.03
$ws.Range("K14").Value = '29/09/2023 02:12'
$ws.Range("L14").Value = 2.78
$ws.Range("M14").Value = '30/09/2023 14:56'
$ws.Range("N14").Value = 3.16
$ws.Range("O14").Value = '29/09/2023 02:12'
$ws.Range("P14").Value = 2.92
$ws.Range("Q14").Value = '30/09/2023 13:25'
$ws.Range("R14").Value = 3.27
$ws.Range("S14").Value = '29/09/2023 02:12'
$ws.Range("T14").Value = 2.72
$ws.Range("U14").Value = '30/09/2023 14:56'
$ws.Range("V14").Value = 'https://www.betexplorer.com/football/greece/super-league-2/panachaiki-panathinaikos/xh8LWHQB/'

$ws.Range("F23").Value = 'AEL Larissa'
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 'Iraklis 1908'
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 1.95
$ws.Range("K23").Value = '10/10/2023 02:12'
$ws.Range("L23").Value = 1.96
$ws.Range("M23").Value = '11/10/2023 14:52'
$ws.Range("N23").Value = 3.06
$ws.Range("O23").Value = '10/10/2023 02:12'
$ws.Range("P23").Value = 3.04
$ws.Range("Q23").Value = '11/10/2023 14:52'
$ws.Range("R23").Value = 3.62
$ws.Range("S23").Value = '10/10/2023 02:12'
$ws.Range("T23").Value = 4.31
$ws.Range("U23").Value = '11/10/2023 14:52'
$ws.Range("V23").Value = 'https://www.betexplorer.com/football/greece/super-league-2/ael-larissa-iraklis-fc/WhhFbXES/'

$ws.Range("F24").Value = 'Karditsa'
$ws.Range("G24").Value = 3
$ws.Range("H24").Value = 'Aiolikos'
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = 1.96
$ws.Range("K24").Value = '10/10/2023 11:12'
$ws.Range("L24").Value = 1.86
$ws.Range("M24").Value = '11/10/2023 14:54'
$ws.Range("N24").Value = 2.97
$ws.Range("O24").Value = '10/10/2023 11:12'
$ws.Range("P24").Value = 3.27
$ws.Range("Q24").Value = '11/10/2023 14:54'
$ws.Range("R24").Value = 3.86
$ws.Range("S24").Value = '10/10/2023 11:12'
$ws.Range("T24").Value = 4.4
$ws.Range("U24").Value = '11/10/2023 14:54'
$ws.Range("V24").Value = 'https://www.betexplorer.com/football/greece/super-league-2/karditsa-aiolikos-fc/2JmBaDaM/'

$ws.Range("F25").Value = 'Kozani FC'
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 'Makedonikos'
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 2.36
$ws.Range("K25").Value = '10/10/2023 11:12'
$ws.Range("L25").Value = 2.63
$ws.Range("M25").Value = '11/10/2023 14:58'
$ws.Range("N25").Value = 2.79
$ws.Range("O25").Value = '10/10/2023 11:12'
$ws.Range("P25").Value = 2.86
$ws.Range("Q25").Value = '11/10/2023 14:58'
$ws.Range("R25").Value = 3.1
$ws.Range("S25").Value = '10/10/2023 11:12'
$ws.Range("T25").Value = 2.93
$ws.Range("U25").Value = '11/10/2023 14:16'
$ws.Range("V25").Value = 'https://www.betexplorer.com/football/greece/super-league-2/kozani-fc-makedonikos-neapolis/t4a21ZU9/'

$ws.Range("F26").Value = 'PAE Egaleo'
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 'Kalamata'
$ws.Range("I26").Value = 3
$ws.Range("J26").Value = 2.62
$ws.Range("K26").Value = '13/10/2023 02:13'
$ws.Range("L26").Value = 4.93
$ws.Range("M26").Value = '14/10/2023 14:52'
$ws.Range("N26").Value = 2.71
$ws.Range("O26").Value = '13/10/2023 02:13'
$ws.Range("P26").Value = 3.53
$ws.Range("Q26").Value = '14/10/2023 14:52'
$ws.Range("R26").Value = 2.75
$ws.Range("S26").Value = '13/10/2023 02:13'
$ws.Range("T26").Value = 1.71
$ws.Range("U26").Value = '14/10/2023 14:52'
$ws.Range("V26").Value = 'https://www.betexplorer.com/football/greece/super-league-2/pae-egaleo-kalamata/dGsbpyeB/'

$ws.Range("F27").Value = 'Giouchtas'
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 'Panachaiki'
$ws.Range("I27").Value = 1
$ws.Range("J27").Value = 1.76
$ws.Range("K27").Value = '14/10/2023 12:13'
$ws.Range("L27").Value = 1.68
$ws.Range("M27").Value = '14/10/2023 14:06'
$ws.Range("N27").Value = 3.44
$ws.Range("O27").Value = '14/10/2023 12:13'
$ws.Range("P27").Value = 3.46
$ws.Range("Q27").Value = '14/10/2023 14:07'
$ws.Range("R27").Value = 4.73
$ws.Range("S27").Value = '14/10/2023 12:13'
$ws.Range("T27").Value = 5.38
$ws.Range("U27").Value = '14/10/2023 14:06'
$ws.Range("V27").Value = 'https://www.betexplorer.com/football/greece/super-league-2/giouchtas-panachaiki/OlUjncQb/'

$ws.Range("F34").Value = 'Athens Kallithea'
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = 'Giouchtas'
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 1.3
$ws.Range("K34").Value = '21/10/2023 01:12'
$ws.Range("L34").Value = 1.23
$ws.Range("M34").Value = '21/10/2023 14:17'
$ws.Range("N34").Value = 4.55
$ws.Range("O34").Value = '21/10/2023 01:12'
$ws.Range("P34").Value = 5.44
$ws.Range("Q34").Value = '21/10/2023 14:18'
$ws.Range("R34").Value = 10.53
$ws.Range("S34").Value = '21/10/2023 01:12'
$ws.Range("T34").Value = 13.52
$ws.Range("U34").Value = '21/10/2023 14:17'
$ws.Range("V34").Value = 'https://www.betexplorer.com/football/greece/super-league-2/athens-kallithea-giouchtas/n5zBsZvU/'

$ws.Range("F36").Value = 'PAOK B'
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 'Aiolikos'
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 1.57
$ws.Range("K36").Value = '21/10/2023 01:12'
$ws.Range("L36").Value = 1.75
$ws.Range("M36").Value = '21/10/2023 14:43'
$ws.Range("N36").Value = 3.91
$ws.Range("O36").Value = '21/10/2023 01:12'
$ws.Range("P36").Value = 3.66
$ws.Range("Q36").Value = '21/10/2023 14:43'
$ws.Range("R36").Value = 5.06
$ws.Range("S36").Value = '21/10/2023 01:12'
$ws.Range("T36").Value = 4.44
$ws.Range("U36").Value = '21/10/2023 14:43'
$ws.Range("V36").Value = 'https://www.betexplorer.com/football/greece/super-league-2/paok-aiolikos-fc/fyiX48y2/'

$ws.Range("F39").Value = 'Panachaiki'
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 'PAE Egaleo'
$ws.Range("I39").Value = 1
$ws.Range("J39").Value = 3.03
$ws.Range("K39").Value = '21/10/2023 02:13'
$ws.Range("L39").Value = 3.1
$ws.Range("M39").Value = '22/10/2023 14:59'
$ws.Range("N39").Value = 2.94
$ws.Range("O39").Value = '21/10/2023 02:13'
$ws.Range("P39").Value = 3.06
$ws.Range("Q39").Value = '22/10/2023 14:57'
$ws.Range("R39").Value = 2.3
$ws.Range("S39").Value = '21/10/2023 02:13'
$ws.Range("T39").Value = 2.37
$ws.Range("U39").Value = '22/10/2023 14:59'
$ws.Range("V39").Value = 'https://www.betexplorer.com/football/greece/super-league-2/panachaiki-pae-egaleo/QV5L1xQo/'

$ws.Range("F40").Value = 'AEL Larissa'
$ws.Range("G40").Value = 1
$ws.Range("H40").Value = 'Makedonikos'
$ws.Range("I40").Value = 1
$ws.Range("J40").Value = 1.67
$ws.Range("K40").Value = '21/10/2023 02:13'
$ws.Range("L40").Value = 1.5
$ws.Range("M40").Value = '22/10/2023 14:56'
$ws.Range("N40").Value = 3.22
$ws.Range("O40").Value = '21/10/2023 02:13'
$ws.Range("P40").Value = 3.72
$ws.Range("Q40").Value = '22/10/2023 14:56'
$ws.Range("R40").Value = 4.89
$ws.Range("S40").Value = '21/10/2023 02:13'
$ws.Range("T40").Value = 7.37
$ws.Range("U40").Value = '22/10/2023 14:56'
$ws.Range("V40").Value = 'https://www.betexplorer.com/football/greece/super-league-2/ael-larissa-makedonikos-neapolis/8KhP6U6k/'

$ws.Range("F41").Value = 'Apollon Pontou'
$ws.Range("G41").Value = 1
$ws.Range("H41").Value = 'Niki Volos'
$ws.Range("I41").Value = 3
$ws.Range("J41").Value = 4.17
$ws.Range("K41").Value = '21/10/2023 02:13'
$ws.Range("L41").Value = 9.99
$ws.Range("M41").Value = '22/10/2023 14:56'
$ws.Range("N41").Value = 3.18
$ws.Range("O41").Value = '21/10/2023 02:13'
$ws.Range("P41").Value = 4.44
$ws.Range("Q41").Value = '22/10/2023 14:56'
$ws.Range("R41").Value = 1.78
$ws.Range("S41").Value = '21/10/2023 02:13'
$ws.Range("T41").Value = 1.34
$ws.Range("U41").Value = '22/10/2023 14:56'
$ws.Range("V41").Value = 'https://www.betexplorer.com/football/greece/super-league-2/apollon-pontou-niki-volos/UTgL7Aiq/'

$ws.Range("F42").Value = 'Kozani FC'
$ws.Range("G42").Value = 1
$ws.Range("H42").Value = 'Karditsa'
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 2.75
$ws.Range("K42").Value = '22/10/2023 12:12'
$ws.Range("L42").Value = 2.58
$ws.Range("M42").Value = '22/10/2023 14:49'
$ws.Range("N42").Value = 2.74
$ws.Range("O42").Value = '22/10/2023 12:12'
$ws.Range("P42").Value = 2.74
$ws.Range("Q42").Value = '22/10/2023 13:54'
$ws.Range("R42").Value = 2.92
$ws.Range("S42").Value = '22/10/2023 12:12'
$ws.Range("T42").Value = 3.14
$ws.Range("U42").Value = '22/10/2023 14:49'
$ws.Range("V42").Value = 'https://www.betexplorer.com/football/greece/super-league-2/kozani-fc-karditsa/OnHhCjqS/'

$ws.Range("F46").Value = 'Kalamata'
$ws.Range("G46").Value = 2
$ws.Range("H46").Value = 'Ilioupoli'
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1.44
$ws.Range("K46").Value = '27/10/2023 02:13'
$ws.Range("L46").Value = 1.38
$ws.Range("M46").Value = '28/10/2023 13:43'
$ws.Range("N46").Value = 3.96
$ws.Range("O46").Value = '27/10/2023 02:13'
$ws.Range("P46").Value = 4.35
$ws.Range("Q46").Value = '28/10/2023 13:43'
$ws.Range("R46").Value = 5.97
$ws.Range("S46").Value = '27/10/2023 02:13'
$ws.Range("T46").Value = 8.56
$ws.Range("U46").Value = '28/10/2023 13:43'
$ws.Range("V46").Value = 'https://www.betexplorer.com/football/greece/super-league-2/kalamata-ilioupoli/bsP9un9q/'

$ws.Range("F47").Value = 'Levadiakos'
$ws.Range("G47").Value = 2
$ws.Range("H47").Value = 'AEK Athens FC B'
$ws.Range("I47").Value = 1
$ws.Range("J47").Value = 1.3
$ws.Range("K47").Value = '28/10/2023 01:12'
$ws.Range("L47").Value = 1.36
$ws.Range("M47").Value = '28/10/2023 13:55'
$ws.Range("N47").Value = 4.78
$ws.Range("O47").Value = '28/10/2023 01:12'
$ws.Range("P47").Value = 4.28
$ws.Range("Q47").Value = '28/10/2023 13:55'
$ws.Range("R47").Value = 9.51
$ws.Range("S47").Value = '28/10/2023 01:12'
$ws.Range("T47").Value = 9.83
$ws.Range("U47").Value = '28/10/2023 13:55'
$ws.Range("V47").Value = 'https://www.betexplorer.com/football/greece/super-league-2/levadiakos-aek/xbk08srf/'

$ws.Range("F48").Value = 'PAE Egaleo'
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = 'Athens Kallithea'
$ws.Range("I48").Value = 1
$ws.Range("J48").Value = 4.46
$ws.Range("K48").Value = '27/10/2023 02:13'
$ws.Range("L48").Value = 7.03
$ws.Range("M48").Value = '28/10/2023 12:23'
$ws.Range("N48").Value = 3.12
$ws.Range("O48").Value = '27/10/2023 02:13'
$ws.Range("P48").Value = 3.72
$ws.Range("Q48").Value = '28/10/2023 12:23'
$ws.Range("R48").Value = 1.75
$ws.Range("S48").Value = '27/10/2023 02:13'
$ws.Range("T48").Value = 1.52
$ws.Range("U48").Value = '28/10/2023 12:23'
$ws.Range("V48").Value = 'https://www.betexplorer.com/football/greece/super-league-2/pae-egaleo-athens-kallithea/WYetcEvH/'

$ws.Range("F49").Value = 'PAE Chania'
$ws.Range("G49").Value = 4
$ws.Range("H49").Value = 'Panachaiki'
$ws.Range("I49").Value = 1
$ws.Range("J49").Value = 1.42
$ws.Range("K49").Value = '27/10/2023 02:13'
$ws.Range("L49").Value = 1.22
$ws.Range("M49").Value = '28/10/2023 13:34'
$ws.Range("N49").Value = 4.07
$ws.Range("O49").Value = '27/10/2023 02:13'
$ws.Range("P49").Value = 5.62
$ws.Range("Q49").Value = '28/10/2023 13:58'
$ws.Range("R49").Value = 6.07
$ws.Range("S49").Value = '27/10/2023 02:13'
$ws.Range("T49").Value = 14.35
$ws.Range("U49").Value = '28/10/2023 13:58'
$ws.Range("V49").Value = 'https://www.betexplorer.com/football/greece/super-league-2/pae-chania-panachaiki/xYNDv6Ok/'

$ws.Range("F50").Value = 'Giouchtas'
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 'Panathinaikos B'
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 1.83
$ws.Range("K50").Value = '28/10/2023 01:12'
$ws.Range("L50").Value = 1.91
$ws.Range("M50").Value = '28/10/2023 13:55'
$ws.Range("N50").Value = 3.35
$ws.Range("O50").Value = '28/10/2023 01:12'
$ws.Range("P50").Value = 3.24
$ws.Range("Q50").Value = '28/10/2023 13:55'
$ws.Range("R50").Value = 4.09
$ws.Range("S50").Value = '28/10/2023 01:12'
$ws.Range("T50").Value = 4.19
$ws.Range("U50").Value = '28/10/2023 13:55'
$ws.Range("V50").Value = 'https://www.betexplorer.com/football/greece/super-league-2/giouchtas-panathinaikos/AJ7XbzA4/'

$ws.Range("F75").Value = 'Levadiakos'
$ws.Range("G75").Value = 2
$ws.Range("H75").Value = 'Aiolikos'
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 1.24
$ws.Range("K75").Value = '11/11/2023 04:43'
$ws.Range("L75").Value = 1.27
$ws.Range("M75").Value = '12/11/2023 11:59'
$ws.Range("N75").Value = 4.72
$ws.Range("O75").Value = '11/11/2023 04:43'
$ws.Range("P75").Value = 4.98
$ws.Range("Q75").Value = '12/11/2023 12:05'
$ws.Range("R75").Value = 10.75
$ws.Range("S75").Value = '11/11/2023 04:43'
$ws.Range("T75").Value = 12.44
$ws.Range("U75").Value = '12/11/2023 11:59'
$ws.Range("V75").Value = 'https://www.betexplorer.com/football/greece/super-league-2/levadiakos-aiolikos-fc/URrrJbdC/'

$ws.Range("F76").Value = 'Niki Volos'
$ws.Range("G76").Value = 3
$ws.Range("H76").Value = 'Iraklis 1908'
$ws.Range("I76").Value = 1
$ws.Range("J76").Value = 1.87
$ws.Range("K76").Value = '11/11/2023 02:12'
$ws.Range("L76").Value = 2.02
$ws.Range("M76").Value = '12/11/2023 12:05'
$ws.Range("N76").Value = 3.03
$ws.Range("O76").Value = '11/11/2023 02:12'
$ws.Range("P76").Value = 3
$ws.Range("Q76").Value = '12/11/2023 12:12'
$ws.Range("R76").Value = 3.98
$ws.Range("S76").Value = '11/11/2023 02:12'
$ws.Range("T76").Value = 4.13
$ws.Range("U76").Value = '12/11/2023 12:05'
$ws.Range("V76").Value = 'https://www.betexplorer.com/football/greece/super-league-2/niki-volos-iraklis-fc/jypvKvs6/'
$ws.Range("A104").Value = 103
$ws.Range("B104").Value = 'greece'
$ws.Range("C104").Value = 'super-league-2'
$ws.Range("D104").Value = '2023-2024'
$ws.Range("E104").Value = 45260.52083333334
$ws.Range("F104").Value = 'AEK Athens FC B'
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 'Kampaniakos'
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 1.91
$ws.Range("K104").Value = '14/11/2023 02:12'
$ws.Range("L104").Value = 1.71
$ws.Range("M104").Value = '30/11/2023 12:23'
$ws.Range("N104").Value = 3.12
$ws.Range("O104").Value = '14/11/2023 02:12'
$ws.Range("P104").Value = 3.5
$ws.Range("Q104").Value = '30/11/2023 12:22'
$ws.Range("R104").Value = 3.71
$ws.Range("S104").Value = '14/11/2023 02:12'
$ws.Range("T104").Value = 4.97
$ws.Range("U104").Value = '30/11/2023 12:23'
$ws.Range("V104").Value = 'https://www.betexplorer.com/football/greece/super-league-2/aek-kampaniakos/MwRNdGJm/'
$ws.Range("A103").Copy()
$ws.Range("A104").PasteSpecial(-4122)
$ws.Range("E103").Copy()
$ws.Range("E104").PasteSpecial(-4122)
$excel.CutCopyMode = 0
